$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 62501748
$ws.Range("I4").Value = 999
$ws.Range("J4").Value = 83335330
$ws.Range("K4").Value = 999
$ws.Range("L4").Value = 83335330
$ws.Range("M4").Value = -885
$ws.Range("N4").Value = -83335558
$ws.Range("H53").Value = 52998.316
$ws.Range("I53").Value = 344.2143
$ws.Range("J53").Value = 200429.8
$ws.Range("K53").Value = 344.2143
$ws.Range("L53").Value = 200429.8
$ws.Range("M53").Value = 292.7857
$ws.Range("N53").Value = -201703.8
$ws.Range("H69").Value = 8552.166999999999
$ws.Range("I69").Value = 5506.5
$ws.Range("J69").Value = 10075
$ws.Range("K69").Value = 16519.5
$ws.Range("L69").Value = 30225
$ws.Range("M69").Value = -15645.5
$ws.Range("N69").Value = -31973
$ws.Range("H72").Value = 8552.166999999999
$ws.Range("I72").Value = 5506.5
$ws.Range("J72").Value = 10075
$ws.Range("K72").Value = 49558.5
$ws.Range("L72").Value = 90675
$ws.Range("M72").Value = -45190.5
$ws.Range("N72").Value = -99411
$ws.Range("H113").Value = 58827748
$ws.Range("I113").Value = 142860240
$ws.Range("K113").Value = 142860240
$ws.Range("M113").Value = -142856986
$ws.Range("H132").Value = 3480.6035
$ws.Range("I132").Value = 3335.125
$ws.Range("K132").Value = 10005.375
$ws.Range("M132").Value = -7475.375
$ws.Range("H138").Value = 1809.6227
$ws.Range("I138").Value = 1426.6061
$ws.Range("J138").Value = 2441.6
$ws.Range("K138").Value = 4279.8183
$ws.Range("L138").Value = 7324.799999999999
$ws.Range("M138").Value = 860.1817000000001
$ws.Range("N138").Value = -17604.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1593.75
$ws.Range("I61").Value = 1593.75
$ws.Range("K61").Value = 1593.75
$ws.Range("M61").Value = -1381.75
$ws.Range("H109").Value = 32188.5
$ws.Range("J109").Value = 32188.5
$ws.Range("L109").Value = 32188.5
$ws.Range("N109").Value = -34962.5
$ws.Range("H122").Value = 2597.244
$ws.Range("I122").Value = 2054.7585
$ws.Range("K122").Value = 6164.2755
$ws.Range("M122").Value = -3714.2755
$ws.Range("H136").Value = 1593.75
$ws.Range("I136").Value = 1593.75
$ws.Range("K136").Value = 4781.25
$ws.Range("M136").Value = -2231.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1183.65
$ws.Range("I107").Value = 1041.6897
$ws.Range("J107").Value = 1557.909
$ws.Range("K107").Value = 1041.6897
$ws.Range("L107").Value = 1557.909
$ws.Range("M107").Value = 878.3103000000001
$ws.Range("N107").Value = -5397.909
$ws.Range("H134").Value = 1126.4348
$ws.Range("I134").Value = 1126.4348
$ws.Range("K134").Value = 3379.3044
$ws.Range("M134").Value = -844.3044

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 720
$ws.Range("I5").Value = 544.875
$ws.Range("J5").Value = 1070.25
$ws.Range("K5").Value = 544.875
$ws.Range("L5").Value = 1070.25
$ws.Range("M5").Value = -432.875
$ws.Range("N5").Value = -1294.25
$ws.Range("H25").Value = 5011
$ws.Range("I25").Value = 5011
$ws.Range("K25").Value = 5011
$ws.Range("M25").Value = -4837
$ws.Range("H122").Value = 3420.2144
$ws.Range("J122").Value = 2879.75
$ws.Range("L122").Value = 8639.25
$ws.Range("N122").Value = -13539.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8990869
$ws.Range("I4").Value = 14175677
$ws.Range("K4").Value = 42527031
$ws.Range("M4").Value = -42526919
$ws.Range("H92").Value = 272.7
$ws.Range("I92").Value = 281.5
$ws.Range("J92").Value = 259.5
$ws.Range("K92").Value = 844.5
$ws.Range("L92").Value = 778.5
$ws.Range("M92").Value = 403.5
$ws.Range("N92").Value = -3274.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4315.5
$ws.Range("I80").Value = 4189.1
$ws.Range("J80").Value = 4947.5
$ws.Range("K80").Value = 4189.1
$ws.Range("L80").Value = 4947.5
$ws.Range("M80").Value = -3191.1
$ws.Range("N80").Value = -6943.5
$ws.Range("H83").Value = 4315.5
$ws.Range("I83").Value = 4189.1
$ws.Range("J83").Value = 4947.5
$ws.Range("K83").Value = 20945.5
$ws.Range("L83").Value = 24737.5
$ws.Range("M83").Value = -15953.5
$ws.Range("N83").Value = -34721.5
$ws.Range("H102").Value = 1517.1111
$ws.Range("I102").Value = 1235
$ws.Range("K102").Value = 1235
$ws.Range("M102").Value = 387
$ws.Range("H113").Value = 1903.95
$ws.Range("I113").Value = 1903.95
$ws.Range("K113").Value = 1903.95
$ws.Range("M113").Value = 266.05
$ws.Range("H122").Value = 3028
$ws.Range("I122").Value = 2299.75
$ws.Range("J122").Value = 3999
$ws.Range("K122").Value = 6899.25
$ws.Range("L122").Value = 11997
$ws.Range("M122").Value = -4449.25
$ws.Range("N122").Value = -16897
$ws.Range("H126").Value = 15114
$ws.Range("I126").Value = 18818.666
$ws.Range("K126").Value = 56455.99800000001
$ws.Range("M126").Value = -53985.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 7000
$ws.Range("I11").Value = 7000
$ws.Range("K11").Value = 7000
$ws.Range("M11").Value = -6860
$ws.Range("H13").Value = 9000
$ws.Range("J13").Value = 9000
$ws.Range("L13").Value = 9000
$ws.Range("N13").Value = -9280
$ws.Range("H16").Value = 690.7646999999999
$ws.Range("J16").Value = 931.8333
$ws.Range("L16").Value = 931.8333
$ws.Range("N16").Value = -1271.8333
$ws.Range("H40").Value = 3249.625
$ws.Range("I40").Value = 2635.9092
$ws.Range("J40").Value = 4599.8
$ws.Range("K40").Value = 2635.9092
$ws.Range("L40").Value = 4599.8
$ws.Range("M40").Value = -2499.9092
$ws.Range("N40").Value = -4871.8
$ws.Range("H132").Value = 5754.0386
$ws.Range("I132").Value = 2932.4
$ws.Range("K132").Value = 8797.200000000001
$ws.Range("M132").Value = -6267.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1756922.1
$ws.Range("I4").Value = 2864443.8
$ws.Range("K4").Value = 2864443.8
$ws.Range("M4").Value = -2864330.8
$ws.Range("H18").Value = 12899
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H113").Value = 1634.7222
$ws.Range("I113").Value = 319.11765
$ws.Range("J113").Value = 24000
$ws.Range("K113").Value = 957.3529500000001
$ws.Range("L113").Value = 72000
$ws.Range("M113").Value = 1212.64705
$ws.Range("N113").Value = -76340
$ws.Range("H122").Value = 1905.1538
$ws.Range("I122").Value = 1508.5
$ws.Range("K122").Value = 4525.5
$ws.Range("M122").Value = -2075.5
$ws.Range("H126").Value = 2063.7856
$ws.Range("I126").Value = 1899.4
$ws.Range("J126").Value = 2474.75
$ws.Range("K126").Value = 5698.200000000001
$ws.Range("L126").Value = 7424.25
$ws.Range("M126").Value = -3228.200000000001
$ws.Range("N126").Value = -12364.25
